# Apply the "Date and Time" / "Cycle_count" update to the analysis sheet.
#
# The original sheet had 42 label/value rows in column A/B. The new layout
# inserts a "Date and Time" row at the very top (shifting everything down
# by one row), renames several labels to include their units, tweaks a
# handful of values, swaps the order of the "highest/lowest cell temp"
# rows, replaces the old "Maximum BMS Temperature in C" row with a new
# "Cycle Count of battery" row, and appends two new speed-bucket rows at
# the bottom. Rather than replay every individual insert/rename, we just
# rewrite the whole A1:B45 range with its final contents directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @(
    "Date and Time",
    "Total time taken for the ride",
    "Actual Ampere-hours (Ah)",
    "Actual Watt-hours (Wh)",
    "Starting SoC (Ah)",
    "Ending SoC (Ah)",
    "Starting SoC (%)",
    "Ending SoC (%)",
    "Total distance covered (km)",
    "Total energy consumption(WH/KM)",
    "Total SOC consumed(%)",
    "Mode",
    "Peak Power(kW)",
    "Average Power(kW)",
    "Total Energy Regenerated(kWh)",
    "Regenerative Effectiveness(%)",
    "Highest Cell Voltage(V)",
    "Lowest Cell Voltage(V)",
    "Difference in Cell Voltage(V)",
    "Minimum Temperature(C)",
    "Maximum Temperature(C)",
    "Difference in Temperature(C)",
    "Maximum Fet Temperature-BMS(C)",
    "Maximum Afe Temperature-BMS(C)",
    "Maximum PCB Temperature-BMS(C)",
    "Maximum MCU Temperature(C)",
    "Maximum Motor Temperature(C)",
    "Abnormal Motor Temperature Detected(C)",
    "highest cell temp(C)",
    "lowest cell temp(C)",
    "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)",
    "Battery Voltage(V)",
    "Total energy charged(kWh)",
    "Electricity consumption units(kW)",
    "Cycle Count of battery",
    "Idling time percentage",
    "Time spent in 0-10 km/h",
    "Time spent in 10-20 km/h",
    "Time spent in 20-30 km/h",
    "Time spent in 30-40 km/h",
    "Time spent in 40-50 km/h",
    "Time spent in 50-60 km/h",
    "Time spent in 60-70 km/h",
    "Time spent in 70-80 km/h",
    "Time spent in 80-90 km/h"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}

# Column B values (row 12 "Mode" keeps its multi-line text value; row 1 is
# also text; everything else is numeric).
# B1 used to hold the old "Total time taken for the ride" duration (with the
# [hh]:mm:ss number format); clear that formatting back to Normal now that
# it holds the new "Date and Time" text.
$ws.Range("B1").Style = "Normal"
$ws.Range("B1").Value = "2024-03-12 20:41:47.242000 to 2024-03-12 21:37:26.915000"

$ws.Range("B2").NumberFormat = "[hh]:mm:ss"
$ws.Range("B2").Value = 0.03843626157407407

$ws.Range("B3").Value = 31.59589166666666
$ws.Range("B4").Value = 1613.489434889166
$ws.Range("B5").Value = 39.532
$ws.Range("B6").Value = 12.329
$ws.Range("B7").Value = 31
$ws.Range("B8").Value = 99
$ws.Range("B9").Value = 32.71985919962416
$ws.Range("B10").Value = 49.31223649360019
$ws.Range("B11").Value = 68

$ws.Range("B12").Value = "Custom mode`n85.53%`nEco mode`n4.29%`nSports mode`n0.22%"

$ws.Range("B13").Value = 6138.107433
$ws.Range("B14").Value = -1756.975791167876
$ws.Range("B15").Value = 72.2175989013889
$ws.Range("B16").Value = 4.284113280289114
$ws.Range("B17").Value = 3.454
$ws.Range("B18").Value = 3.062
$ws.Range("B19").Value = 0.3920000000000003
$ws.Range("B20").Value = 34
$ws.Range("B21").Value = 44
$ws.Range("B22").Value = 10
$ws.Range("B23").Value = 70
$ws.Range("B24").Value = 63
$ws.Range("B25").Value = 61
$ws.Range("B26").Value = 48
$ws.Range("B27").Value = 0
$ws.Range("B28").Value = 1
$ws.Range("B29").Value = 44
$ws.Range("B30").Value = -1
$ws.Range("B31").Value = 45
$ws.Range("B32").Value = 55
$ws.Range("B33").Value = 1.737774041666666
$ws.Range("B34").Value = [double]"1.453960878235163e-07"
$ws.Range("B35").Value = 42
$ws.Range("B36").Value = 17.61544284632854
$ws.Range("B37").Value = 9.367903103709311
$ws.Range("B38").Value = 3.917486752460257
$ws.Range("B39").Value = 8.096139288417865
$ws.Range("B40").Value = 11.65404996214989
$ws.Range("B41").Value = 16.27933383800151
$ws.Range("B42").Value = 15.96139288417865
$ws.Range("B43").Value = 16.06358819076457
$ws.Range("B44").Value = 0.9348978046934141
$ws.Range("B45").Value = 0
